# Fix the CDS "read table" Cypher query in cell B3 ("startup" sheet):
# the Tumor column should read the sample's tumor-status property directly
# instead of the unresolved `tumor` alias, and tidy the ORDER BY indent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$queryLines = @(
    'MATCH (s:study)<--(p:participant)<--(samp:sample)',
    'WHERE s.study_name in ["GECCO OICR: Molecular Pathological Epidemiology of Colorectal Cancer"]',
    'WITH p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor',
    'RETURN  ',
    ' coalesce(samp.sample_id, '''') as `Sample ID`,',
    ' coalesce(p.participant_id,'''') as `Participant ID`,',
    ' coalesce(s.study_name, '''') as `Study Name`,',
    ' coalesce(s.phs_accession,'''') as `Accession`,',
    ' coalesce(samp.sample_tumor_status,'''') as `Tumor`,',
    'coalesce(samp.sample_type,'''') as `Analyte Type`',
    '  ORDER By samp.sample_id LIMIT 100'
)
$newQuery = $queryLines -join "`r`n"

$ws.Range("B3").Value = $newQuery

# Move the selection onto the cell that was just edited (matches the
# author re-selecting B3 after fixing the query).
$ws.Range("B3").Select()
